# Update the runs/balls/fours/sixes figures for Rahul Tewatia's innings.
# Source data keeps these as text-typed numerals (the sheet stores
# numberStoredAsText for A1:F12), so values are entered with a leading
# apostrophe to force text, and ClearFormats() strips the resulting
# quote-prefix style so the cell formatting matches the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> new (runs, balls, fours, sixes) values, columns C:F
$newValues = @{
    2  = @("24", "12", "0", "3")
    3  = @("19", "11", "1", "1")
    4  = @("2",  "3",  "0", "0")
    5  = @("14", "18", "1", "0")
    6  = @("53", "31", "0", "7")
    8  = @("10", "8",  "1", "0")
    9  = @("31", "27", "2", "1")
    10 = @("38", "29", "3", "2")
    11 = @("5",  "6",  "1", "0")
    12 = @("14", "10", "0", "1")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $cell = $ws.Cells.Item($row, 3 + $i)
        $cell.Value = "'" + $vals[$i]
        $cell.ClearFormats()
    }
}
